$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update translated text in shared strings (via cell values) ---

# New Kyrgyz title for A1 (replaces the old, now-unused Kyrgyz title)
$ws.Range("A1").Value = "4.2.2.1а Балдарды мектепке чейин билим берүү менен камтуу"

# "urban settlements" -> "urban" for the Kyrgyz "urban" rows (column C)
foreach ($r in 6,9,12,15,18,21,24,27) {
    $ws.Cells.Item($r, 3).Value = "urban"
}

# "countryside" -> "rural" for the Kyrgyz "rural" rows (column C)
foreach ($r in 7,10,13,16,19,22,25,28) {
    $ws.Cells.Item($r, 3).Value = "rural"
}

# --- Add new column N (year 2023) ---

# Header row: year label 2023, formatted like the other year header cells
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# Data rows 5-30: value plus formatting copied from a same-styled neighbor cell
$rowInfo = @{
    5  = @{ Value = 28.34784779265912;  Source = "K5"  }
    6  = @{ Value = 39.999446500300472; Source = "K6"  }
    7  = @{ Value = 23.198557483143556; Source = "K7"  }
    8  = @{ Value = 27.597876990321573; Source = "K8"  }
    9  = @{ Value = 47.175678010018999; Source = "K9"  }
    10 = @{ Value = 22.17579894112394;  Source = "K10" }
    11 = @{ Value = 24.100104034215697; Source = "K11" }
    12 = @{ Value = 38.296287676015361; Source = "K12" }
    13 = @{ Value = 19.410249509822766; Source = "K13" }
    14 = @{ Value = 30.400174646089773; Source = "K14" }
    15 = @{ Value = 44.562134629854725; Source = "K15" }
    16 = @{ Value = 24.612036336109007; Source = "K16" }
    17 = @{ Value = 39.266683582846994; Source = "K17" }
    18 = @{ Value = 54.818496110630946; Source = "K18" }
    19 = @{ Value = 36.591078066914498; Source = "K19" }
    20 = @{ Value = 23.890520476423561; Source = "K20" }
    21 = @{ Value = 16.93085228577992;  Source = "K21" }
    22 = @{ Value = 24.386979772654026; Source = "K22" }
    23 = @{ Value = 28.919699950811605; Source = "K23" }
    24 = @{ Value = 37.932834522359492; Source = "K24" }
    25 = @{ Value = 26.985549456704376; Source = "K25" }
    26 = @{ Value = 27.190143693828379; Source = "K26" }
    27 = @{ Value = 54.006768771869439; Source = "K27" }
    28 = @{ Value = 22.334624692306893; Source = "K28" }
    29 = @{ Value = 36.01461582008131;  Source = "K29" }
    30 = @{ Value = 42.081208505725009; Source = "M30" }
}

foreach ($r in 5..30) {
    $info = $rowInfo[$r]
    $ws.Range($info.Source).Copy()
    $ws.Range("N$r").PasteSpecial(-4122)
    $ws.Range("N$r").Value = $info.Value
    $ws.Range("N$r").NumberFormat = "0.0"
}

$excel.CutCopyMode = 0
